# Purchase orders template: add a status parameter to the report header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row layout -------------------------------------------------
# Previously A1:G1 was one big merged title cell ("Approved purchase
# orders"). Split it so A1:C1 keeps a (shorter) title, and free up
# D1/E1/F1/G1 to host a new "Status:" label + {header:title} placeholder.
$ws.Range("A1:G1").UnMerge()
$ws.Range("A1:C1").Merge()

$ws.Rows(1).RowHeight = 23

$ws.Range("A1").Value = "הזמנות רכש "
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = ""
$ws.Range("F1").Value = "{header:title}"
$ws.Range("E1").Value = "סטטוס:"
$ws.Range("G1").Value = ""

# Vertically center everything in the (now taller) header row.
$ws.Range("A1:I1").VerticalAlignment = -4108

# E1 ("Status:") uses the plain (non-underlined) label font, like the
# rest of the workbook's field labels.
$ws.Range("E1").Font.Underline = -4142

# H1/I1 swap their horizontal alignment (date label goes to the left of
# the placeholder, which now sits on the right).
$ws.Range("H1").HorizontalAlignment = -4131
$ws.Range("I1").HorizontalAlignment = -4152

# Give the new A1:C1 title a thin bottom rule to visually separate it
# from the status area.
$bottom = $ws.Range("A1:C1").Borders.Item(9)
$bottom.ColorIndex = 15
$bottom.Weight = 2

# Reflect the new selection left behind by the editor.
$ws.Range("E5").Select()
